$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.065.63"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "'1.818.98"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'310.59"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.5010"
$ws.Range("E7").Value = "  -2.38%  "
$ws.Range("D8").Value = "'0.3914"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").Value = "'0.09975"
$ws.Range("E9").Value = "  +27.35%  "
$ws.Range("D10").Value = "'1.109"
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("D11").Value = "'40.81"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "'6.429"
$ws.Range("D13").Value = "'20.57"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").Value = "'1.000"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "'1.816.93"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").Value = "'7.297"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("E17").Value = "  +6.06%  "
$ws.Range("D18").Value = "'92.49"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").Value = "'0.06651"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "'0.9992"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "'17.20"
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("D22").Value = "'5.954"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").Value = "'28.115.40"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").Value = "'11.17"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").Value = "'2.248"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").Value = "'159.00"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("D27").Value = "'20.74"
$ws.Range("E27").Value = "  +2.17%  "
$ws.Range("D28").Value = "'2.027.44"
$ws.Range("E28").Value = "  +1.81%  "
$ws.Range("E29").Value = "  +3.01%  "
$ws.Range("D30").Value = "'127.01"
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("D31").Value = "'0.1067"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "'1.038"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").Value = "'3.616"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").Value = "'0.06718"
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("D37").Value = "'8.915"
$ws.Range("E37").Value = "  +2.24%  "
$ws.Range("D38").Value = "'0.2140"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").Value = "'4.958"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").Value = "'11.34"
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("D41").Value = "'0.6195"
$ws.Range("E41").Value = "  +1.58%  "
$ws.Range("D42").Value = "'1.174"
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("D43").Value = "'0.9988"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "'13.25"
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").Value = "'0.5913"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.283"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.689"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").Value = "'124.48"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "'1.940"
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("D50").Value = "'1.182"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("D51").Value = "'0.06787"
$ws.Range("E51").Value = "  -0.44%  "
